# Weekly update: insert the newest week's Betarraga (beet) prices for
# Terminal Hortofrutícola Agro Chillán at the top of the data block
# (row 488), pushing the existing historical rows down by two rows.
#
# This mirrors the source diff: dimension grows from A1:R554 to A1:R556,
# and two brand-new rows (Primera / Segunda quality) are inserted right
# after the existing row 487.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 488-489; everything that used to start at row
# 488 shifts down to row 490 onward (dimension becomes A1:R556).
$ws.Range("488:489").Insert()

# New row 488: "Primera" quality, newest reporting week.
$ws.Cells.Item(488, 1).Value  = 7
$ws.Cells.Item(488, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(488, 3).Value  = "Ñuble"
$ws.Cells.Item(488, 4).Value  = 44984
$ws.Cells.Item(488, 5).Value  = 16
$ws.Cells.Item(488, 6).Value  = 100114014
$ws.Cells.Item(488, 7).Value  = "Betarraga"
$ws.Cells.Item(488, 8).Value  = "Sin especificar"
$ws.Cells.Item(488, 9).Value  = "Primera"
$ws.Cells.Item(488, 10).Value = 400
$ws.Cells.Item(488, 11).Value = 700
$ws.Cells.Item(488, 12).Value = 800
$ws.Cells.Item(488, 13).Value = 750
$ws.Cells.Item(488, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(488, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(488, 16).Value = 150
$ws.Cells.Item(488, 17).Value = 5
$ws.Cells.Item(488, 18).Value = "Hortaliza"

# New row 489: "Segunda" quality, same reporting week.
$ws.Cells.Item(489, 1).Value  = 7
$ws.Cells.Item(489, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(489, 3).Value  = "Ñuble"
$ws.Cells.Item(489, 4).Value  = 44984
$ws.Cells.Item(489, 5).Value  = 16
$ws.Cells.Item(489, 6).Value  = 100114014
$ws.Cells.Item(489, 7).Value  = "Betarraga"
$ws.Cells.Item(489, 8).Value  = "Sin especificar"
$ws.Cells.Item(489, 9).Value  = "Segunda"
$ws.Cells.Item(489, 10).Value = 300
$ws.Cells.Item(489, 11).Value = 600
$ws.Cells.Item(489, 12).Value = 600
$ws.Cells.Item(489, 13).Value = 600
$ws.Cells.Item(489, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(489, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(489, 16).Value = 120
$ws.Cells.Item(489, 17).Value = 5
$ws.Cells.Item(489, 18).Value = "Hortaliza"

# Date column keeps its datetime number format after the insert, but make
# sure explicitly in case the host didn't copy row formatting down.
$ws.Range("D488:D489").NumberFormat = "YYYY-MM-DD HH:MM:SS"
